$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (original row 26) and the "SC 92" row
# (original row 28, which becomes row 27 after the first deletion).
# Remaining rows shift up, filling the gap left by the two removed
# records, and the sheet's used range shrinks from F35 to F33.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Fill in / clear individual cells that were (re)imputed, now that
# the rows below have shifted up into their final positions.
$ws.Range("C3").Value = 11.2
$ws.Range("E4").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("E9").Value = -6.8
$ws.Range("E10").Value = -6.1
$ws.Range("E11").Value = -7.9
$ws.Range("E12").Value = -5.3
$ws.Range("F12").Value = ""
$ws.Range("F13").Value = 17.1
$ws.Range("F14").Value = 17.76
$ws.Range("E15").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = 17.78
$ws.Range("E18").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("F25").Value = ""

# The remaining SC rows (now rows 27-33 after the deletions above)
# also had some of their own missing values imputed.
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("E31").Value = -8.1
$ws.Range("C32").Value = 10.5
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39
